$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-06-22 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-23 Friday", 2) | Out-Null

# Update each math-problem cell in the practice table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "79-37="
$t.Cell(1,2).Range.Text = "14+85="
$t.Cell(1,3).Range.Text = "33+0="
$t.Cell(1,4).Range.Text = "85-57="
$t.Cell(1,5).Range.Text = "75-20="
$t.Cell(2,1).Range.Text = "84-57="
$t.Cell(2,2).Range.Text = "43-28="
$t.Cell(2,3).Range.Text = "15+80="
$t.Cell(2,4).Range.Text = "73-57="
$t.Cell(2,5).Range.Text = "77-60="
$t.Cell(3,1).Range.Text = "30-4="
$t.Cell(3,2).Range.Text = "28+39="
$t.Cell(3,3).Range.Text = "1-0="
$t.Cell(3,4).Range.Text = "54-5="
$t.Cell(3,5).Range.Text = "78-40="
$t.Cell(4,1).Range.Text = "93-75="
$t.Cell(4,2).Range.Text = "64-24="
$t.Cell(4,3).Range.Text = "70-60="
$t.Cell(4,4).Range.Text = "21-7="
$t.Cell(4,5).Range.Text = "49+29="
$t.Cell(5,1).Range.Text = "7+20="
$t.Cell(5,2).Range.Text = "72-42="
$t.Cell(5,3).Range.Text = "68+23="
$t.Cell(5,4).Range.Text = "36+63="
$t.Cell(5,5).Range.Text = "94-93="
$t.Cell(6,1).Range.Text = "18+72="
$t.Cell(6,2).Range.Text = "56+2="
$t.Cell(6,3).Range.Text = "99-38="
$t.Cell(6,4).Range.Text = "63-42="
$t.Cell(6,5).Range.Text = "35+25="
$t.Cell(7,1).Range.Text = "61+21="
$t.Cell(7,2).Range.Text = "15+23="
$t.Cell(7,3).Range.Text = "27+37="
$t.Cell(7,4).Range.Text = "27-26="
$t.Cell(7,5).Range.Text = "8+65="
$t.Cell(8,1).Range.Text = "86-76="
$t.Cell(8,2).Range.Text = "49-20="
$t.Cell(8,3).Range.Text = "0+43="
$t.Cell(8,4).Range.Text = "83-4="
$t.Cell(8,5).Range.Text = "2+3="
$t.Cell(9,1).Range.Text = "80-77="
$t.Cell(9,2).Range.Text = "9+4="
$t.Cell(9,3).Range.Text = "21-9="
$t.Cell(9,4).Range.Text = "58+11="
$t.Cell(9,5).Range.Text = "70+15="
$t.Cell(10,1).Range.Text = "74+24="
$t.Cell(10,2).Range.Text = "82-15="
$t.Cell(10,3).Range.Text = "0+3="
$t.Cell(10,4).Range.Text = "46-40="
$t.Cell(10,5).Range.Text = "23+59="
$t.Cell(11,1).Range.Text = "74-62="
$t.Cell(11,2).Range.Text = "40+6="
$t.Cell(11,3).Range.Text = "89+7="
$t.Cell(11,4).Range.Text = "24+66="
$t.Cell(11,5).Range.Text = "8+67="
$t.Cell(12,1).Range.Text = "78-39="
$t.Cell(12,2).Range.Text = "67-0="
$t.Cell(12,3).Range.Text = "82-11="
$t.Cell(12,4).Range.Text = "31+23="
$t.Cell(12,5).Range.Text = "67+11="
$t.Cell(13,1).Range.Text = "17+47="
$t.Cell(13,2).Range.Text = "92-6="
$t.Cell(13,3).Range.Text = "70-21="
$t.Cell(13,4).Range.Text = "31+46="
$t.Cell(13,5).Range.Text = "3-1="
$t.Cell(14,1).Range.Text = "68+26="
$t.Cell(14,2).Range.Text = "11+25="
$t.Cell(14,3).Range.Text = "67-16="
$t.Cell(14,4).Range.Text = "39-11="
$t.Cell(14,5).Range.Text = "41+34="
$t.Cell(15,1).Range.Text = "20+26="
$t.Cell(15,2).Range.Text = "60-37="
$t.Cell(15,3).Range.Text = "44+11="
$t.Cell(15,4).Range.Text = "37+17="
$t.Cell(15,5).Range.Text = "15+10="
$t.Cell(16,1).Range.Text = "75-26="
$t.Cell(16,2).Range.Text = "28+55="
$t.Cell(16,3).Range.Text = "77-69="
$t.Cell(16,4).Range.Text = "93-44="
$t.Cell(16,5).Range.Text = "67-20="
$t.Cell(17,1).Range.Text = "33-15="
$t.Cell(17,2).Range.Text = "3+65="
$t.Cell(17,3).Range.Text = "15+18="
$t.Cell(17,4).Range.Text = "38+6="
$t.Cell(17,5).Range.Text = "12-0="
$t.Cell(18,1).Range.Text = "45+35="
$t.Cell(18,2).Range.Text = "48-37="
$t.Cell(18,3).Range.Text = "6+37="
$t.Cell(18,4).Range.Text = "91-29="
$t.Cell(18,5).Range.Text = "6+9="
$t.Cell(19,1).Range.Text = "44+5="
$t.Cell(19,2).Range.Text = "36-19="
$t.Cell(19,3).Range.Text = "68+14="
$t.Cell(19,4).Range.Text = "32+7="
$t.Cell(19,5).Range.Text = "58+23="
$t.Cell(20,1).Range.Text = "49-17="
$t.Cell(20,2).Range.Text = "73-51="
$t.Cell(20,3).Range.Text = "26-0="
$t.Cell(20,4).Range.Text = "45+19="
$t.Cell(20,5).Range.Text = "73-41="

Write-Host "Done updating date and" $t.Rows.Count "x" $t.Columns.Count "table."
